$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns before column B, shifting old B:V to K:AE
$ws.Range("B:J").Insert()

# New header values for the 9 newly inserted weeks (row 1, columns B..J)
$ws.Range("B1").Value = "Sep_08"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("J1").Value = "Jun_16"

# Fill new columns B..J with "UN" for all analyst rows (2..33)
$ws.Range("B2:J33").Value = "UN"

# New annotation in C5 (Zacks Investment Research, week Aug_25): downgrade event
# Copy the existing "Downgrades" highlight style (orange fill) from U5 (an existing downgrade annotation)
$ws.Range("U5").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C5").Value = "8/22/2019,Downgrades,Strong-Buy -> Hold,"
